$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number cells (e.g. "604.79") to stay as text instead of
# being auto-converted to a numeric value by Excels input parser.
$textCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D38",
    "D39",
    "D41",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = '67.755.72'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '3.498.81'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '604.79'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("D6").Value = '149.41'
$ws.Range("E6").Value = '  -1.72%  '
$ws.Range("D7").Value = '3.496.67'
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.485'
$ws.Range("E9").Value = '  +0.64%  '
$ws.Range("D10").Value = '0.143'
$ws.Range("E10").Value = '  +2.92%  '
$ws.Range("D11").Value = '7.54'
$ws.Range("E11").Value = '  +6.16%  '
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").Value = '0.0000214'
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("D14").Value = '31.86'
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("D15").Value = '4.090.65'
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.515.61'
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '67.833.50'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("E19").Value = '  +1.42%  '
$ws.Range("D20").Value = '15.32'
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("D21").Value = '9.94'
$ws.Range("E21").Value = '  +2.35%  '
$ws.Range("D22").Value = '444.64'
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("D23").Value = '0.623'
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("D24").Value = '78.91'
$ws.Range("E24").Value = '  +2.50%  '
$ws.Range("D25").Value = '3.641.15'
$ws.Range("E25").Value = '  -1.13%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '5.68'
$ws.Range("E27").Value = '  -4.10%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0000123'
$ws.Range("E28").Value = '  -4.19%  '
$ws.Range("D29").Value = '9.96'
$ws.Range("E29").Value = '  -2.03%  '
$ws.Range("D30").Value = '8.63'
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("D31").Value = '2.50'
$ws.Range("E31").Value = '  -1.49%  '
$ws.Range("D32").Value = '1.64'
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("D33").Value = '0.169'
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = '25.52'
$ws.Range("E35").Value = '  -0.55%  '
$ws.Range("D36").Value = '6.15'
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("B37").Value = 'RenzoRestakedETH'
$ws.Range("C37").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D37").Value = '3.494.70'
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '1.84'
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").Value = '7.96'
$ws.Range("E39").Value = '  -0.85%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").Value = '2.32'
$ws.Range("E41").Value = '  +6.36%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '176.71'
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("D44").Value = '0.0899'
$ws.Range("E44").Value = '  +0.65%  '
$ws.Range("D45").Value = '5.39'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").Value = '0.895'
$ws.Range("E46").Value = '  +0.96%  '
$ws.Range("D47").Value = '30.01'
$ws.Range("E47").Value = '  +4.98%  '
$ws.Range("D48").Value = '46.72'
$ws.Range("E48").Value = '  +2.70%  '
$ws.Range("D49").Value = '1.28'
$ws.Range("E49").Value = '  +0.50%  '
$ws.Range("D50").Value = '2.51'
$ws.Range("E50").Value = '  -5.79%  '
$ws.Range("D51").Value = '7.59'
$ws.Range("E51").Value = '  +0.05%  '

# Restore default style on the forced-text cells so no stray number format
# style lingers on them (matches original workbook styling).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
